# GCP/txt_ctr_table.xlsx — add "no vote intervention" baseline table
#
# 1. Rename existing "Sheet1" -> "baseline_upvote"
# 2. Add a new sheet "baseline_novote" right after it, with a small
#    treatment/control table (upvote-only -> baseline vs no-vote -> baseline)
# 3. Leave the selection on baseline_upvote at A1:G3 and activate baseline_novote

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "baseline_upvote"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "baseline_novote"

# ---- formatting first (copy existing styles so we reuse style indexes) ----

# Header row A1:G1 -> same formatting as baseline_upvote's header row
$ws1.Range("A1:G1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122)

# Row 2 numeric cells (baseline n / no votes n columns)
$ws1.Range("E9:F9").Copy()
$ws2.Range("E2:F2").PasteSpecial(-4122)

# Row 2 date cells (baseline dates / no votes dates columns)
$ws1.Range("C18").Copy()
$ws2.Range("C2:D2").PasteSpecial(-4122)

# Row 3 numeric + correlation cells
$ws1.Range("E3:G3").Copy()
$ws2.Range("E3:G3").PasteSpecial(-4122)

# Row 3 date cells
$ws1.Range("C19").Copy()
$ws2.Range("C3:D3").PasteSpecial(-4122)

# ---- values (write order controls new shared-string index assignment) ----

$ws2.Range("A1").Value = "treatment subreddit"
$ws2.Range("B1").Value = "control subreddit"
$ws2.Range("C1").Value = "baseline dates"
$ws2.Range("E1").Value = "baseline n"
$ws2.Range("G1").Value = "correlation"

$ws2.Range("A2").Value = "\unpopularopinion"
$ws2.Range("D2").Value = "02-21-2018 to 04-01-2018"
$ws2.Range("C2").Value = "01-01-2010 to 02-20-2018"
$ws2.Range("F1").Value = "no votes n"
$ws2.Range("D1").Value = "no votes dates"

# "no votes n" header: italicize+bold the trailing "n" like the other n-headers
$nChar = $ws2.Range("F1").Characters(10, 1)
$nChar.Font.Bold = $true
$nChar.Font.Italic = $true

$ws2.Range("E2").Value = 612169
$ws2.Range("F2").Value = 151280

$ws2.Range("B3").Value = "\PoliticalHumor"
$ws2.Range("E3").Value = 1519118
$ws2.Range("F3").Value = 210184
$ws2.Range("G3").Value = 0.979

# ---- cosmetics: column widths, selections, active sheet ----

$ws2.Range("A1:F3").EntireColumn.AutoFit() | Out-Null

$ws1.Range("A1:G3").Select()
$ws2.Range("G4").Select()
$ws2.Activate()
